# Add "2022-Q4" data to the 600022-山东钢铁 fund-holdings workbook.
#
# Shape of the edit:
#   1) A brand-new quarter sheet "2022-Q4" is inserted right after the
#      "总计" (summary) sheet and right before the existing "2022-Q3"
#      sheet. It reuses the same fund list/formatting as "2022-Q3" (the
#      most recently-known quarter) but carries fresh position numbers.
#   2) The "总计" summary sheet gets a new top data row for 2022-Q4 and
#      every older quarter row shifts down by one.
#   3) All the other existing quarter sheets are untouched (they simply
#      move one tab to the right).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet by duplicating "2022-Q3" (keeps the
#    fund code/name/size columns + header formatting identical).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# Fresh position data for the new quarter. Columns D-G store their
# numbers as plain text in this workbook (matching the rest of the file),
# so format those cells as Text before typing the values in; column H is
# a genuine number.
$textCells = $q4.Range("E2:G5,D5")
$textCells.NumberFormat = "@"

$q4.Cells.Item(2, 5).Value2 = "92.42"
$q4.Cells.Item(2, 6).Value2 = "1.43"
$q4.Cells.Item(2, 7).Value2 = "0.0122"
$q4.Cells.Item(2, 8).Value2 = 9

$q4.Cells.Item(3, 5).Value2 = "91.91"
$q4.Cells.Item(3, 6).Value2 = "3.13"
$q4.Cells.Item(3, 7).Value2 = "0.0025"

$q4.Cells.Item(4, 5).Value2 = "92.42"
$q4.Cells.Item(4, 6).Value2 = "1.43"
$q4.Cells.Item(4, 7).Value2 = "0.0006"
$q4.Cells.Item(4, 8).Value2 = 9

$q4.Cells.Item(5, 4).Value2 = "0.02"
$q4.Cells.Item(5, 5).Value2 = "91.91"
$q4.Cells.Item(5, 6).Value2 = "3.13"
$q4.Cells.Item(5, 7).Value2 = "0.0006"

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q4 and
#    push the older quarters down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Copy formatting from the row below (keeps column A's style) before
# writing the new values.
$summary.Cells.Item(3, 1).Copy($summary.Cells.Item(2, 1))

$summary.Cells.Item(2, 1).Value2 = 0
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 4
$summary.Cells.Item(2, 4).Value2 = 0.02

$summary.Cells.Item(3, 1).Value2 = 1
$summary.Cells.Item(3, 2).Value2 = "2022-Q3"
$summary.Cells.Item(3, 3).Value2 = 4
$summary.Cells.Item(3, 4).Value2 = 0.02

$summary.Cells.Item(4, 1).Value2 = 2
$summary.Cells.Item(4, 2).Value2 = "2021-Q4"
$summary.Cells.Item(4, 3).Value2 = 2
$summary.Cells.Item(4, 4).Value2 = 0.28

$summary.Cells.Item(5, 1).Value2 = 3
$summary.Cells.Item(5, 2).Value2 = "2021-Q1"
$summary.Cells.Item(5, 3).Value2 = 1
$summary.Cells.Item(5, 4).Value2 = 0.08

$summary.Cells.Item(6, 1).Value2 = 4
$summary.Cells.Item(6, 2).Value2 = "2020-Q4"
$summary.Cells.Item(6, 3).Value2 = 1
$summary.Cells.Item(6, 4).Value2 = 0.03

# ---------------------------------------------------------------------
# 3) Restore the original active tab ("2020-Q4" is the last sheet and
#    was the selected one before this edit).
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
